$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plans are changing dynamically: split the double/triple lesson labels
# into individual lesson numbers.
$ws.Range("A4").Value = "Урок№1"
$ws.Range("A5").Value = "Урок№2"
$ws.Range("A6").Value = "Урок№3"
$ws.Range("A7").Value = "Урок№4"
$ws.Range("A8").Value = "Урок№5"
$ws.Range("A9").Value = "Урок№6"
$ws.Range("A10").Value = "Урок№7"

# Move the active selection off the table, matching the author's last cursor position.
$ws.Range("D17").Select()
